$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Tornado paragraph
#    a) "...fulfill the strike role.  Italy converted 16 of these..."
#       -> "...fulfill the strike role.  Eventually converting 16 of these..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Italy converted 16 of these",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "Eventually converting 16 of these", 2) | Out-Null

#    b) "...more than a decade away (in 2003), a solutions was needed."
#       -> "...more than a decade away (in 2003), a quick solutions was needed."
$d.Content.Find.Execute(
    "(in 2003), a solutions was needed",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "(in 2003), a quick solutions was needed", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. F-104 Starfighter paragraph
#    a) "...for a time but it became the backbone..."
#       -> "...for a time but as it's popularity faded in many countries, it became the backbone..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "for a time but it became the backbone",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "for a time but as it’s popularity faded in many countries, it became the backbone", 2) | Out-Null

#    b) "...was the ultimate version with improved radar..."
#       -> "...was the ultimate fighter version with improved radar..."
$d.Content.Find.Execute(
    "was the ultimate version with improved radar",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "was the ultimate fighter version with improved radar", 2) | Out-Null

#    c) "147 of these were subsequently modified" -> "147 were subsequently modified"
$d.Content.Find.Execute(
    "147 of these were subsequently modified",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "147 were subsequently modified", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. A-11 squadron table: "18" -> "16" (both the Indians/Treviso and Black
#    Cats/Villafranca rows, column 4 "No.")
# ---------------------------------------------------------------------------
$amxTable = $d.Tables.Item(3)
$amxTable.Cell(2, 4).Range.Text = "16"
$amxTable.Cell(3, 4).Range.Text = "16"

# ---------------------------------------------------------------------------
# 4. G.91Y / AMX paragraph
#    a) "By 1994 the subsonic aircraft" -> "By 1994 this subsonic aircraft"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "By 1994 the subsonic aircraft",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "By 1994 this subsonic aircraft", 2) | Out-Null

#    b) "...procured to replace it. 65 aircraft were operated."
#       -> "...procured to replace it. In Northern Fury, 65 aircraft remained with 48 employed in operational squadrons."
$d.Content.Find.Execute(
    "procured to replace it. 65 aircraft were operated.",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "procured to replace it. In Northern Fury, 65 aircraft remained with 48 employed in operational squadrons.", 2) | Out-Null
